$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.091.51"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "2.620.52"
$ws.Range("E3").Value = "  +4.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.10"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.76"
$ws.Range("E6").Value = "  +2.38%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "2.618.19"
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("E10").Value = "  +15.20%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +2.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.03"
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("D14").Value = "3.089.71"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.63"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000183"
$ws.Range("E16").Value = "  +7.61%  "
$ws.Range("D17").Value = "71.057.48"
$ws.Range("E17").Value = "  +4.26%  "
$ws.Range("D18").Value = "2.624.03"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "381.83"
$ws.Range("E19").Value = "  +8.68%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.95"
$ws.Range("E20").Value = "  +5.63%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.55"
$ws.Range("E21").Value = "  +4.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.49"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.46"
$ws.Range("E24").Value = "  +5.65%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.86"
$ws.Range("E26").Value = "  +6.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  +4.32%  "
$ws.Range("D28").Value = "2.747.68"
$ws.Range("E28").Value = "  +4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.988"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").Value = "0.0₃0958"
$ws.Range("E30").Value = "  +6.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "534.22"
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("E32").Value = "  +2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.32"
$ws.Range("E33").Value = "  +3.64%  "
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "164.06"
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("E37").Value = "  -1.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.19"
$ws.Range("E38").Value = "  +4.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.89"
$ws.Range("E39").Value = "  +6.56%  "
$ws.Range("E40").Value = "  +1.44%  "
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +9.75%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.04"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.332"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.08"
$ws.Range("E46").Value = "  +2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.81"
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("E48").Value = "  +1.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.532"
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0266"
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.68"
$ws.Range("E51").Value = "  +5.27%  "
